# "Last run before holiday / Changes data in and out"
#
# The "_set_needs" sheet listed "BEV battery charge" as its own need,
# separate from "BEV battery discharge" a row below it - a leftover/duplicate
# entry. Remove that entire row; Excel then drops the now-unused
# "BEV battery charge" shared string and renumbers every other shared-string
# reference across the workbook accordingly.

$wb = $excel.ActiveWorkbook

$wsNeeds = $wb.Worksheets.Item("_set_needs")
$wsNeeds.Activate()
$wsNeeds.Rows("8:8").Delete()
$wsNeeds.Rows("8:8").Select()

# With that row gone, the technology names on "_set_t" (e.g. "Air
# conditioning unit") no longer fit column A - widen it to fit.
$wsT = $wb.Worksheets.Item("_set_t")
$wsT.Columns("A:A").ColumnWidth = 17
